# MeleeWeaponDB_Sheet.xlsx edit
# Renames header columns to camelCase field names, fixes typo/path strings,
# and corrects a couple of numeric stat values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 : header renames (PascalCase -> camelCase field names) ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "itemName"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "icon"
$ws.Range("E1").Value = "iconWidth"
$ws.Range("F1").Value = "iconHeight"
$ws.Range("G1").Value = "itemPrefab"
$ws.Range("H1").Value = "price"
$ws.Range("I1").Value = "atk"
$ws.Range("J1").Value = "atkRate"
$ws.Range("K1").Value = "critRate"
$ws.Range("L1").Value = "critDamage"
$ws.Range("M1").Value = "range"
$ws.Range("N1").Value = "lifeSteal"
$ws.Range("O1").Value = "type"
$ws.Range("P1").Value = "weaponTier"

# --- Row 2 : Axe ---
$ws.Range("C2").Value = "발등 조심하세요."
$ws.Range("D2").Value = "Resources/Icons/Weapons"
$ws.Range("G2").Value = "Resources/Prefabs/Weapons"
$ws.Range("K2").Value = 1

# --- Row 3 : Spear ---
$ws.Range("C3").Value = "창"
$ws.Range("D3").Value = "Resources/Icons/Weapons"
$ws.Range("G3").Value = "Resources/Prefabs/Weapons"
$ws.Range("K3").Value = 1
$ws.Range("N3").Value = 1

# --- Sheet view: clear frozen/scrolled top-left cell, move selection ---
$ws.Range("O3").Select()
